# Applies the "added pom and DD concepta for automation" edit:
#  - Rename Sheet1 -> login
#  - Update A1's phone-code label from "(+91)India" to "(+1)India"
#  - Move the active selection to A6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "login"

# Fix the India calling-code text in A1 (was mistakenly "(+91)India").
$ws.Range("A1").Value = "(+1)India"

# Leave the selection on A6, matching the saved workbook view state.
$ws.Range("A6").Select()
